# Initial Glider Design v5.0 - results refresh
# Updates the cached numeric results in column B (Sheet1) to reflect the
# re-run of the Python optimizer after the roll-in/tangent-point constraint
# change and plotting refinements described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.153149525267364
$ws.Range("B5").Value = 1.262260630994263
$ws.Range("B7").Value = 0.08119501541474206
$ws.Range("B8").Value = 8.789208926400274
$ws.Range("B9").Value = 1.387313746705807
$ws.Range("B10").Value = 1.386477061009285
$ws.Range("B11").Value = -(0.0008366856965216218)
$ws.Range("B12").Value = 0.04000007475393155
$ws.Range("B13").Value = 0.4000000353888141
$ws.Range("B14").Value = 0.0199999926655473
$ws.Range("B15").Value = 49656.65013336878
$ws.Range("B16").Value = -(0.09342439026005291)
$ws.Range("B17").Value = 0.6968318429726539
$ws.Range("B18").Value = 0.7948883238098
$ws.Range("B19").Value = 2.03539828202274
$ws.Range("B20").Value = 0.1407574956027135
$ws.Range("B22").Value = 0.1118864897433062
$ws.Range("B23").Value = 0.6968318429726539
$ws.Range("B24").Value = 0.1901606294779517
$ws.Range("B25").Value = 0.04754015736948792
$ws.Range("B27").Value = 0.009040266250862706
$ws.Range("B28").Value = 0.6968318429726539
$ws.Range("B29").Value = 0.07145091670240844
$ws.Range("B30").Value = 0.03572545835120422
$ws.Range("B32").Value = 0.002552616748807255
$ws.Range("B36").Value = -(8.350203523124514 * [Math]::Pow(10, -15))
$ws.Range("B37").Value = -(9.848384207410734 * [Math]::Pow(10, -15))
$ws.Range("B38").Value = -(3.026541962318201)
$ws.Range("B39").Value = 1.058550505592157
$ws.Range("B40").Value = 60.65047637187003
$ws.Range("B41").Value = 1.071714576406072
$ws.Range("B42").Value = -(31.62631796216882)
$ws.Range("B43").Value = 0.01065729830535309
$ws.Range("B45").Value = 7.000429547638715 * [Math]::Pow(10, -7)
$ws.Range("B47").Value = 0.0105649986999016
$ws.Range("B48").Value = 9.159956249672909 * [Math]::Pow(10, -5)
$ws.Range("B49").Value = 8.971452785922185 * [Math]::Pow(10, -12)
$ws.Range("B50").Value = 0.0267508736520472
$ws.Range("B51").Value = -(1.488107865180076 * [Math]::Pow(10, -18))
$ws.Range("B52").Value = 0.003113164800574826
$ws.Range("B53").Value = 0.00183568148991003
$ws.Range("B54").Value = 0.001454769286758737
$ws.Range("B55").Value = 0.0032867990891481
$ws.Range("B56").Value = 3.353689553173407 * [Math]::Pow(10, -21)
$ws.Range("B57").Value = -(6.2297072224905 * [Math]::Pow(10, -6))
$ws.Range("B58").Value = 2.2815618898642 * [Math]::Pow(10, -21)
$ws.Range("B59").Value = 0.03232607677269975
$ws.Range("B60").Value = 0.000777269098820298
$ws.Range("B61").Value = 0.0001649274672767937
$ws.Range("B66").Value = 0.007112306099094361
$ws.Range("B68").Value = -(9.609426332622078 * [Math]::Pow(10, -9))
$ws.Range("B69").Value = 0.006014445586277189
$ws.Range("B70").Value = -(0.03518937390067838)
$ws.Range("B73").Value = -(0.03518937390067838)
$ws.Range("B74").Value = 0.397193404367839
$ws.Range("B75").Value = 0.01411599656328559
$ws.Range("B76").Value = 0.649291685603166
$ws.Range("B79").Value = 0.649291685603166
$ws.Range("B80").Value = 0.09508031473897584
$ws.Range("B82").Value = 0.6611063846214497
$ws.Range("B85").Value = 0.6611063846214497
$ws.Range("B87").Value = 0.07145091670240844
$ws.Range("B88").Value = -(0.1161011467611064)
$ws.Range("B89").Value = 6.415188180862725 * [Math]::Pow(10, -17)
$ws.Range("B90").Value = 1.893337490654271
$ws.Range("B91").Value = -(3.821464854012173 * [Math]::Pow(10, -16))
$ws.Range("B92").Value = 4.904675092714389 * [Math]::Pow(10, -11)
$ws.Range("B93").Value = -(6.002786140652568 * [Math]::Pow(10, -17))
$ws.Range("B94").Value = 0.1161011467611064
$ws.Range("B95").Value = -(0.2144372955744961)
$ws.Range("B96").Value = 3.821464854012173 * [Math]::Pow(10, -16)
$ws.Range("B97").Value = 3.867645503516725 * [Math]::Pow(10, -16)
$ws.Range("B98").Value = 1.884734192416254
$ws.Range("B99").Value = 6.415188180862725 * [Math]::Pow(10, -17)
$ws.Range("B100").Value = 0.2144372955744961
$ws.Range("B101").Value = 3.821464854012173 * [Math]::Pow(10, -16)
$ws.Range("B102").Value = 4.904675092714389 * [Math]::Pow(10, -11)
$ws.Range("B103").Value = 6.002786140652568 * [Math]::Pow(10, -17)
$ws.Range("B104").Value = 1.035669480373438
$ws.Range("B105").Value = 3.525173277221761 * [Math]::Pow(10, -17)
$ws.Range("B106").Value = 0.1178342088629518
$ws.Range("B107").Value = 2.641768886918701 * [Math]::Pow(10, -16)
$ws.Range("B108").Value = 1.914739871437279 * [Math]::Pow(10, -10)
$ws.Range("B109").Value = 4.149710717489231 * [Math]::Pow(10, -17)
$ws.Range("B110").Value = 1.820765797462427
$ws.Range("B111").Value = 0.06728617070302481
$ws.Range("B112").Value = 4.683753385137379 * [Math]::Pow(10, -17)
$ws.Range("B113").Value = 3.816391647148976 * [Math]::Pow(10, -16)
$ws.Range("B114").Value = 0.04024134293396982
$ws.Range("B115").Value = 7.112366251504909 * [Math]::Pow(10, -17)
$ws.Range("B116").Value = 0.7943869669449269
$ws.Range("B117").Value = 0.8006279918397183
$ws.Range("B118").Value = 0.06228093452468567
$ws.Range("B119").Value = 0.007013516317176
$ws.Range("B121").Value = -(5.421010862427522 * [Math]::Pow(10, -20))
$ws.Range("B122").Value = -(0.0396493032369044)
$ws.Range("B123").Value = -(0)
$ws.Range("B124").Value = 0.1901606294779517
$ws.Range("B126").Value = -(2.710505431213761 * [Math]::Pow(10, -20))
$ws.Range("B127").Value = 0.001214686381204971
$ws.Range("B128").Value = 1.731434795725346 * [Math]::Pow(10, -17)
$ws.Range("B129").Value = 5.643467968323278 * [Math]::Pow(10, -19)
$ws.Range("B130").Value = -(9.620023459473031 * [Math]::Pow(10, -5))
$ws.Range("B131").Value = -(1.109731046703648 * [Math]::Pow(10, -17))
$ws.Range("B132").Value = 0.07145091670240844
$ws.Range("B134").Value = 0.001687460429141284
$ws.Range("B135").Value = 0.001327472025569427
$ws.Range("B137").Value = -(2.816001888311887 * [Math]::Pow(10, -21))
$ws.Range("B138").Value = -(0.0004958394134239406)
$ws.Range("B139").Value = 1.509358513064274 * [Math]::Pow(10, -21)
$ws.Range("B142").Value = 0.07684184542697522
$ws.Range("B143").Value = 0.1375954501475208
$ws.Range("B144").Value = 4.064860030510662
$ws.Range("B145").Value = 0.6767864455493355
$ws.Range("B146").Value = -(1.4746458108678 * [Math]::Pow(10, -12))
$ws.Range("B147").Value = -(1.511397943680217 * [Math]::Pow(10, -11))
$ws.Range("B148").Value = -(0.162594705084695)
$ws.Range("B149").Value = -(2.81709631988402 * [Math]::Pow(10, -12))
$ws.Range("B150").Value = 0.03238118399833192
$ws.Range("B151").Value = -(5.506125375674737 * [Math]::Pow(10, -6))
$ws.Range("B152").Value = -(6.063441250288863 * [Math]::Pow(10, -7))
$ws.Range("B153").Value = 0.1356588869903509
$ws.Range("B154").Value = -(0.02500007760707386)
$ws.Range("B155").Value = 7.541842719580456 * [Math]::Pow(10, -7)
$ws.Range("B156").Value = 0.02026731492435541
$ws.Range("B157").Value = -(0.09200472242985269)
$ws.Range("B158").Value = -(0.007946217166265869)
$ws.Range("B159").Value = -(0.00172887961528978)
$ws.Range("B160").Value = -(0.07824093433272834)
$ws.Range("B161").Value = -(0.6483472533641985)
$ws.Range("B162").Value = -(0.000311782902418352)
$ws.Range("B163").Value = -(0.273078453147326)
$ws.Range("B164").Value = 0.3851218934731725
$ws.Range("B165").Value = 0.5315074849041035
$ws.Range("B166").Value = -(1.499129589323552 * [Math]::Pow(10, -14))
$ws.Range("B167").Value = -(1.58294944266581 * [Math]::Pow(10, -13))
$ws.Range("B168").Value = -(10.53213465556191)
$ws.Range("B169").Value = -(3.119970347510272 * [Math]::Pow(10, -14))
$ws.Range("B170").Value = -(0.0001383543122290831)
$ws.Range("B171").Value = -(3.183313593280168 * [Math]::Pow(10, -5))
$ws.Range("B172").Value = 0.08964986524589011
$ws.Range("B173").Value = 0.2217406558473109
$ws.Range("B174").Value = -(9.254732585931308 * [Math]::Pow(10, -5))
$ws.Range("B175").Value = -(0.04708434991402697)

Write-Output "Updated 144 result cells in column B"
